$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last three data rows (old rows 5-7, the "ECs"/"FAPs" sending-cluster
# combinations are no longer part of the recomputed (new-TPM) output).
$ws.Range("A5:A7").EntireRow.Delete()

# Row 2: MuSCs -> Fgf15 -> Fgfr1 -> ECs, with recomputed TPM-based statistics.
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Fgf15"
$ws.Range("C2").Value = "Fgfr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.0005903333333333333
$ws.Range("H2").Value = 0.001771
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.845768666666667
$ws.Range("N2").Value = 5.537306
$ws.Range("O2").Value = 0.01459089321241885
$ws.Range("P2").Value = 0.01459089321241885
$ws.Range("Q2").Value = 0.001089618769555555
$ws.Range("R2").Value = 0.009806568926000001
$ws.Range("S2").Value = 0.01459089321241885
$ws.Range("T2").Value = 0.01459089321241885

# Row 3: MuSCs -> Fgf15 -> Fgfr1 -> FAPs, with recomputed TPM-based statistics.
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Fgf15"
$ws.Range("C3").Value = "Fgfr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.0005903333333333333
$ws.Range("H3").Value = 0.001771
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 82.95722966666666
$ws.Range("N3").Value = 248.871689
$ws.Range("O3").Value = 0.6557810310272387
$ws.Range("P3").Value = 0.6557810310272387
$ws.Range("Q3").Value = 0.04897241791322222
$ws.Range("R3").Value = 0.440751761219
$ws.Range("S3").Value = 0.6557810310272387
$ws.Range("T3").Value = 0.6557810310272387

# Row 4: MuSCs -> Fgf15 -> Fgfr1 -> MuSCs, with recomputed TPM-based statistics.
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Fgf15"
$ws.Range("C4").Value = "Fgfr1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.0005903333333333333
$ws.Range("H4").Value = 0.001771
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 41.69841866666667
$ws.Range("N4").Value = 125.095256
$ws.Range("O4").Value = 0.3296280757603424
$ws.Range("P4").Value = 0.3296280757603424
$ws.Range("Q4").Value = 0.02461596648622223
$ws.Range("R4").Value = 0.221543698376
$ws.Range("S4").Value = 0.3296280757603424
$ws.Range("T4").Value = 0.3296280757603424
